$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Efnb1"
$ws.Cells.Item(2,3).Value = "Ephb1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3.0
$ws.Cells.Item(2,6).Value = 1.0
$ws.Cells.Item(2,7).Value = 9.546140333333334
$ws.Cells.Item(2,8).Value = 28.638421
$ws.Cells.Item(2,9).Value = 0.587227294878132
$ws.Cells.Item(2,10).Value = 0.587227294878132
$ws.Cells.Item(2,11).Value = 3.0
$ws.Cells.Item(2,12).Value = 1.0
$ws.Cells.Item(2,13).Value = 4.002008666666667
$ws.Cells.Item(2,14).Value = 12.006026
$ws.Cells.Item(2,15).Value = 0.4834231243738785
$ws.Cells.Item(2,16).Value = 0.4834231243738787
$ws.Cells.Item(2,17).Value = 38.20373634721622
$ws.Cells.Item(2,18).Value = 343.833627124946
$ws.Cells.Item(2,19).Value = 0.2838792536076074
$ws.Cells.Item(2,20).Value = 0.2838792536076075

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Efnb1"
$ws.Cells.Item(3,3).Value = "Ephb1"
$ws.Cells.Item(3,4).Value = "MuSCs"
$ws.Cells.Item(3,5).Value = 3.0
$ws.Cells.Item(3,6).Value = 1.0
$ws.Cells.Item(3,7).Value = 9.546140333333334
$ws.Cells.Item(3,8).Value = 28.638421
$ws.Cells.Item(3,9).Value = 0.587227294878132
$ws.Cells.Item(3,10).Value = 0.587227294878132
$ws.Cells.Item(3,11).Value = 3.0
$ws.Cells.Item(3,12).Value = 1.0
$ws.Cells.Item(3,13).Value = 4.265473333333333
$ws.Cells.Item(3,14).Value = 12.79642
$ws.Cells.Item(3,15).Value = 0.5152483708764571
$ws.Cells.Item(3,16).Value = 0.5152483708764573
$ws.Cells.Item(3,17).Value = 40.71880702809111
$ws.Cells.Item(3,18).Value = 366.46926325282
$ws.Cells.Item(3,19).Value = 0.3025679070201464
$ws.Cells.Item(3,20).Value = 0.3025679070201465

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Efnb1"
$ws.Cells.Item(4,3).Value = "Ephb1"
$ws.Cells.Item(4,4).Value = "Resolving-Mac"
$ws.Cells.Item(4,5).Value = 3.0
$ws.Cells.Item(4,6).Value = 1.0
$ws.Cells.Item(4,7).Value = 9.546140333333334
$ws.Cells.Item(4,8).Value = 28.638421
$ws.Cells.Item(4,9).Value = 0.587227294878132
$ws.Cells.Item(4,10).Value = 0.587227294878132
$ws.Cells.Item(4,11).Value = 1.0
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 0.010998
$ws.Cells.Item(4,14).Value = 0.032994
$ws.Cells.Item(4,15).Value = 0.001328504749664189
$ws.Cells.Item(4,16).Value = 0.00132850474966419
$ws.Cells.Item(4,17).Value = 0.104988451386
$ws.Cells.Item(4,18).Value = 0.9448960624740002
$ws.Cells.Item(4,19).Value = 0.0007801342503780518
$ws.Cells.Item(4,20).Value = 0.0007801342503780521

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Efnb1"
$ws.Cells.Item(5,3).Value = "Ephb1"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3.0
$ws.Cells.Item(5,6).Value = 1.0
$ws.Cells.Item(5,7).Value = 4.058683666666667
$ws.Cells.Item(5,8).Value = 12.176051
$ws.Cells.Item(5,9).Value = 0.2496684258894083
$ws.Cells.Item(5,10).Value = 0.2496684258894083
$ws.Cells.Item(5,11).Value = 3.0
$ws.Cells.Item(5,12).Value = 1.0
$ws.Cells.Item(5,13).Value = 4.002008666666667
$ws.Cells.Item(5,14).Value = 12.006026
$ws.Cells.Item(5,15).Value = 0.4834231243738785
$ws.Cells.Item(5,16).Value = 0.4834231243738787
$ws.Cells.Item(5,17).Value = 16.24288720925844
$ws.Cells.Item(5,18).Value = 146.185984883326
$ws.Cells.Item(5,19).Value = 0.1206954905009659
$ws.Cells.Item(5,20).Value = 0.1206954905009659

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Efnb1"
$ws.Cells.Item(6,3).Value = "Ephb1"
$ws.Cells.Item(6,4).Value = "MuSCs"
$ws.Cells.Item(6,5).Value = 3.0
$ws.Cells.Item(6,6).Value = 1.0
$ws.Cells.Item(6,7).Value = 4.058683666666667
$ws.Cells.Item(6,8).Value = 12.176051
$ws.Cells.Item(6,9).Value = 0.2496684258894083
$ws.Cells.Item(6,10).Value = 0.2496684258894083
$ws.Cells.Item(6,11).Value = 3.0
$ws.Cells.Item(6,12).Value = 1.0
$ws.Cells.Item(6,13).Value = 4.265473333333333
$ws.Cells.Item(6,14).Value = 12.79642
$ws.Cells.Item(6,15).Value = 0.5152483708764571
$ws.Cells.Item(6,16).Value = 0.5152483708764573
$ws.Cells.Item(6,17).Value = 17.31220694860222
$ws.Cells.Item(6,18).Value = 155.80986253742
$ws.Cells.Item(6,19).Value = 0.1286412496988071
$ws.Cells.Item(6,20).Value = 0.1286412496988071

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Efnb1"
$ws.Cells.Item(7,3).Value = "Ephb1"
$ws.Cells.Item(7,4).Value = "Resolving-Mac"
$ws.Cells.Item(7,5).Value = 3.0
$ws.Cells.Item(7,6).Value = 1.0
$ws.Cells.Item(7,7).Value = 4.058683666666667
$ws.Cells.Item(7,8).Value = 12.176051
$ws.Cells.Item(7,9).Value = 0.2496684258894083
$ws.Cells.Item(7,10).Value = 0.2496684258894083
$ws.Cells.Item(7,11).Value = 1.0
$ws.Cells.Item(7,12).Value = 0.3333333333333333
$ws.Cells.Item(7,13).Value = 0.010998
$ws.Cells.Item(7,14).Value = 0.032994
$ws.Cells.Item(7,15).Value = 0.001328504749664189
$ws.Cells.Item(7,16).Value = 0.00132850474966419
$ws.Cells.Item(7,17).Value = 0.04463740296600001
$ws.Cells.Item(7,18).Value = 0.4017366266940001
$ws.Cells.Item(7,19).Value = 0.0003316856896352605
$ws.Cells.Item(7,20).Value = 0.0003316856896352606

# Row 8
$ws.Cells.Item(8,1).Value = "MuSCs"
$ws.Cells.Item(8,2).Value = "Efnb1"
$ws.Cells.Item(8,3).Value = "Ephb1"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3.0
$ws.Cells.Item(8,6).Value = 1.0
$ws.Cells.Item(8,7).Value = 2.210442
$ws.Cells.Item(8,8).Value = 6.631326
$ws.Cells.Item(8,9).Value = 0.1359745227725727
$ws.Cells.Item(8,10).Value = 0.1359745227725727
$ws.Cells.Item(8,11).Value = 3.0
$ws.Cells.Item(8,12).Value = 1.0
$ws.Cells.Item(8,13).Value = 4.002008666666667
$ws.Cells.Item(8,14).Value = 12.006026
$ws.Cells.Item(8,15).Value = 0.4834231243738785
$ws.Cells.Item(8,16).Value = 0.4834231243738787
$ws.Cells.Item(8,17).Value = 8.846208041164001
$ws.Cells.Item(8,18).Value = 79.61587237047601
$ws.Cells.Item(8,19).Value = 0.06573322863396416
$ws.Cells.Item(8,20).Value = 0.06573322863396419

# Row 9
$ws.Cells.Item(9,1).Value = "MuSCs"
$ws.Cells.Item(9,2).Value = "Efnb1"
$ws.Cells.Item(9,3).Value = "Ephb1"
$ws.Cells.Item(9,4).Value = "MuSCs"
$ws.Cells.Item(9,5).Value = 3.0
$ws.Cells.Item(9,6).Value = 1.0
$ws.Cells.Item(9,7).Value = 2.210442
$ws.Cells.Item(9,8).Value = 6.631326
$ws.Cells.Item(9,9).Value = 0.1359745227725727
$ws.Cells.Item(9,10).Value = 0.1359745227725727
$ws.Cells.Item(9,11).Value = 3.0
$ws.Cells.Item(9,12).Value = 1.0
$ws.Cells.Item(9,13).Value = 4.265473333333333
$ws.Cells.Item(9,14).Value = 12.79642
$ws.Cells.Item(9,15).Value = 0.5152483708764571
$ws.Cells.Item(9,16).Value = 0.5152483708764573
$ws.Cells.Item(9,17).Value = 9.428581405880001
$ws.Cells.Item(9,18).Value = 84.85723265292
$ws.Cells.Item(9,19).Value = 0.07006065133927178
$ws.Cells.Item(9,20).Value = 0.0700606513392718

# Row 10
$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Efnb1"
$ws.Cells.Item(10,3).Value = "Ephb1"
$ws.Cells.Item(10,4).Value = "Resolving-Mac"
$ws.Cells.Item(10,5).Value = 3.0
$ws.Cells.Item(10,6).Value = 1.0
$ws.Cells.Item(10,7).Value = 2.210442
$ws.Cells.Item(10,8).Value = 6.631326
$ws.Cells.Item(10,9).Value = 0.1359745227725727
$ws.Cells.Item(10,10).Value = 0.1359745227725727
$ws.Cells.Item(10,11).Value = 1.0
$ws.Cells.Item(10,12).Value = 0.3333333333333333
$ws.Cells.Item(10,13).Value = 0.010998
$ws.Cells.Item(10,14).Value = 0.032994
$ws.Cells.Item(10,15).Value = 0.001328504749664189
$ws.Cells.Item(10,16).Value = 0.00132850474966419
$ws.Cells.Item(10,17).Value = 0.024310441116
$ws.Cells.Item(10,18).Value = 0.218793970044
$ws.Cells.Item(10,19).Value = 0.0001806427993366843
$ws.Cells.Item(10,20).Value = 0.0001806427993366843

# Row 11
$ws.Cells.Item(11,1).Value = "Resolving-Mac"
$ws.Cells.Item(11,2).Value = "Efnb1"
$ws.Cells.Item(11,3).Value = "Ephb1"
$ws.Cells.Item(11,4).Value = "ECs"
$ws.Cells.Item(11,5).Value = 3.0
$ws.Cells.Item(11,6).Value = 1.0
$ws.Cells.Item(11,7).Value = 0.4410293333333333
$ws.Cells.Item(11,8).Value = 1.323088
$ws.Cells.Item(11,9).Value = 0.02712975645988715
$ws.Cells.Item(11,10).Value = 0.02712975645988715
$ws.Cells.Item(11,11).Value = 3.0
$ws.Cells.Item(11,12).Value = 1.0
$ws.Cells.Item(11,13).Value = 4.002008666666667
$ws.Cells.Item(11,14).Value = 12.006026
$ws.Cells.Item(11,15).Value = 0.4834231243738785
$ws.Cells.Item(11,16).Value = 0.4834231243738787
$ws.Cells.Item(11,17).Value = 1.765003214254222
$ws.Cells.Item(11,18).Value = 15.885028928288
$ws.Cells.Item(11,19).Value = 0.01311515163134106
$ws.Cells.Item(11,20).Value = 0.01311515163134106

# Row 12
$ws.Cells.Item(12,1).Value = "Resolving-Mac"
$ws.Cells.Item(12,2).Value = "Efnb1"
$ws.Cells.Item(12,3).Value = "Ephb1"
$ws.Cells.Item(12,4).Value = "MuSCs"
$ws.Cells.Item(12,5).Value = 3.0
$ws.Cells.Item(12,6).Value = 1.0
$ws.Cells.Item(12,7).Value = 0.4410293333333333
$ws.Cells.Item(12,8).Value = 1.323088
$ws.Cells.Item(12,9).Value = 0.02712975645988715
$ws.Cells.Item(12,10).Value = 0.02712975645988715
$ws.Cells.Item(12,11).Value = 3.0
$ws.Cells.Item(12,12).Value = 1.0
$ws.Cells.Item(12,13).Value = 4.265473333333333
$ws.Cells.Item(12,14).Value = 12.79642
$ws.Cells.Item(12,15).Value = 0.5152483708764571
$ws.Cells.Item(12,16).Value = 0.5152483708764573
$ws.Cells.Item(12,17).Value = 1.881198860551111
$ws.Cells.Item(12,18).Value = 16.93078974496
$ws.Cells.Item(12,19).Value = 0.01397856281823189
$ws.Cells.Item(12,20).Value = 0.0139785628182319

# Row 13
$ws.Cells.Item(13,1).Value = "Resolving-Mac"
$ws.Cells.Item(13,2).Value = "Efnb1"
$ws.Cells.Item(13,3).Value = "Ephb1"
$ws.Cells.Item(13,4).Value = "Resolving-Mac"
$ws.Cells.Item(13,5).Value = 3.0
$ws.Cells.Item(13,6).Value = 1.0
$ws.Cells.Item(13,7).Value = 0.4410293333333333
$ws.Cells.Item(13,8).Value = 1.323088
$ws.Cells.Item(13,9).Value = 0.02712975645988715
$ws.Cells.Item(13,10).Value = 0.02712975645988715
$ws.Cells.Item(13,11).Value = 1.0
$ws.Cells.Item(13,12).Value = 0.3333333333333333
$ws.Cells.Item(13,13).Value = 0.010998
$ws.Cells.Item(13,14).Value = 0.032994
$ws.Cells.Item(13,15).Value = 0.001328504749664189
$ws.Cells.Item(13,16).Value = 0.00132850474966419
$ws.Cells.Item(13,17).Value = 0.004850440608
$ws.Cells.Item(13,18).Value = 0.043653965472
$ws.Cells.Item(13,19).Value = 0.0000360420103141928
$ws.Cells.Item(13,20).Value = 0.00003604201031419282

